# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
# (GitHub Actions cryptos-list refresh: Mon Oct  7 10:10:46 UTC 2024)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a literal string into a cell without Excel's automatic
# type-sniffing turning numeric-looking / boolean-looking text into a
# real number/bool (e.g. '0.0000180' -> 1.8E-05, '1.00' -> 1).
# Temporarily mark the cell as Text, assign, then restore the 'Normal'
# style so no stray number-format override is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# --- Row 30 / 31 swapped places (Coin name + Link) ---
Set-TextValue $ws.Range("B30") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("B31") "Fetch.AI"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"

# --- Price (column D) and Volume(1h) (column E) refresh ---
Set-TextValue $ws.Range("D2") "63.168.49"
Set-TextValue $ws.Range("E2") "  +1.78%  "
Set-TextValue $ws.Range("D3") "2.449.26"
Set-TextValue $ws.Range("E3") "  +1.14%  "
Set-TextValue $ws.Range("E4") "  +0.09%  "
Set-TextValue $ws.Range("D5") "572.17"
Set-TextValue $ws.Range("E5") "  +1.59%  "
Set-TextValue $ws.Range("D6") "147.10"
Set-TextValue $ws.Range("E6") "  +2.33%  "
Set-TextValue $ws.Range("E7") "  -0.06%  "
Set-TextValue $ws.Range("E8") "  +0.91%  "
Set-TextValue $ws.Range("D9") "2.452.52"
Set-TextValue $ws.Range("E9") "  +1.32%  "
Set-TextValue $ws.Range("E10") "  +2.68%  "
Set-TextValue $ws.Range("E11") "  +0.45%  "
Set-TextValue $ws.Range("E12") "  +2.31%  "
Set-TextValue $ws.Range("D13") "5.25"
Set-TextValue $ws.Range("E13") "  +0.96%  "
Set-TextValue $ws.Range("D14") "26.92"
Set-TextValue $ws.Range("E14") "  +2.66%  "
Set-TextValue $ws.Range("D15") "0.0000180"
Set-TextValue $ws.Range("E15") "  +3.72%  "
Set-TextValue $ws.Range("D16") "2.884.10"
Set-TextValue $ws.Range("E16") "  +0.86%  "
Set-TextValue $ws.Range("D17") "63.204.02"
Set-TextValue $ws.Range("E17") "  +1.97%  "
Set-TextValue $ws.Range("D18") "2.451.19"
Set-TextValue $ws.Range("E18") "  +1.10%  "
Set-TextValue $ws.Range("D19") "11.37"
Set-TextValue $ws.Range("E19") "  +1.12%  "
Set-TextValue $ws.Range("D20") "7.24"
Set-TextValue $ws.Range("E20") "  +5.96%  "
Set-TextValue $ws.Range("D21") "326.33"
Set-TextValue $ws.Range("E21") "  +0.94%  "
Set-TextValue $ws.Range("D22") "4.19"
Set-TextValue $ws.Range("E22") "  +1.30%  "
Set-TextValue $ws.Range("E23") "  +12.31%  "
Set-TextValue $ws.Range("D24") "0.994"
Set-TextValue $ws.Range("E24") "  -0.45%  "
Set-TextValue $ws.Range("D25") "66.66"
Set-TextValue $ws.Range("E25") "  -0.94%  "
Set-TextValue $ws.Range("D26") "618.89"
Set-TextValue $ws.Range("E26") "  +11.04%  "
Set-TextValue $ws.Range("E27") "  +0.74%  "
Set-TextValue $ws.Range("D28") "0.0000104"
Set-TextValue $ws.Range("E28") "  +10.94%  "
Set-TextValue $ws.Range("D29") "2.590.85"
Set-TextValue $ws.Range("E29") "  +1.93%  "
Set-TextValue $ws.Range("D30") "0.999"
Set-TextValue $ws.Range("E30") "  -0.01%  "
Set-TextValue $ws.Range("D31") "1.49"
Set-TextValue $ws.Range("E31") "  +6.76%  "
Set-TextValue $ws.Range("D32") "8.25"
Set-TextValue $ws.Range("E32") "  +0.29%  "
Set-TextValue $ws.Range("D33") "0.144"
Set-TextValue $ws.Range("E33") "  -2.13%  "
Set-TextValue $ws.Range("E34") "  +1.99%  "
Set-TextValue $ws.Range("D35") "5.13"
Set-TextValue $ws.Range("E35") "  +8.13%  "
Set-TextValue $ws.Range("E36") "  -0.03%  "
Set-TextValue $ws.Range("E37") "  -0.05%  "
Set-TextValue $ws.Range("D38") "0.382"
Set-TextValue $ws.Range("E38") "  +0.50%  "
Set-TextValue $ws.Range("D39") "18.76"
Set-TextValue $ws.Range("E39") "  +0.45%  "
Set-TextValue $ws.Range("E40") "  -1.51%  "
Set-TextValue $ws.Range("D41") "1.80"
Set-TextValue $ws.Range("E41") "  +0.19%  "
Set-TextValue $ws.Range("D42") "145.73"
Set-TextValue $ws.Range("E42") "  -4.16%  "
Set-TextValue $ws.Range("E43") "  +16.52%  "
Set-TextValue $ws.Range("E44") "  -0.12%  "
Set-TextValue $ws.Range("D45") "148.24"
Set-TextValue $ws.Range("E45") "  +0.52%  "
Set-TextValue $ws.Range("E47") "  +1.80%  "
Set-TextValue $ws.Range("D48") "20.74"
Set-TextValue $ws.Range("E48") "  +4.02%  "
Set-TextValue $ws.Range("D49") "0.601"
Set-TextValue $ws.Range("E49") "  +1.18%  "
Set-TextValue $ws.Range("E50") "  +3.32%  "
Set-TextValue $ws.Range("D51") "0.0925"
Set-TextValue $ws.Range("E51") "  +0.45%  "
